$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Intake metric: mean_Intake (O) and sem_Intake (P) were only the raw milligrams
# intaken; now divided by each animals body weight to give a proper intake rate.
$ws.Range("O2").Value = 662.60393509215135
$ws.Range("P2").Value = 192.7298011765437
$ws.Range("O3").Value = 954.82091884984868
$ws.Range("P3").Value = 344.05848332991911
$ws.Range("O4").Value = 812.07754296646465
$ws.Range("P4").Value = 283.08670479648754
$ws.Range("O5").Value = 955.27029418837765
$ws.Range("P5").Value = 279.45777242579493
$ws.Range("O6").Value = 899.30636025498472
$ws.Range("P6").Value = 219.30690905617593
$ws.Range("O7").Value = 1038.636728649863
$ws.Range("P7").Value = 243.52309346660232
$ws.Range("O8").Value = 1048.1940092526545
$ws.Range("P8").Value = 217.04241147763923
$ws.Range("O9").Value = 1286.0533069705627
$ws.Range("P9").Value = 369.79571776552183
$ws.Range("O10").Value = 943.99066899691809
$ws.Range("P10").Value = 181.66081318784734
$ws.Range("O11").Value = 810.31570844193584
$ws.Range("P11").Value = 169.81642932644624
$ws.Range("O12").Value = 805.07085469088599
$ws.Range("P12").Value = 160.90806161211731
$ws.Range("O13").Value = 925.89455499344774
$ws.Range("P13").Value = 196.74481614021289
$ws.Range("O14").Value = 885.88373992487504
$ws.Range("P14").Value = 216.88905013462497
$ws.Range("O15").Value = 713.65236589105518
$ws.Range("P15").Value = 168.38894238238058
$ws.Range("O16").Value = 924.25520616131644
$ws.Range("P16").Value = 202.33334121095055
$ws.Range("O17").Value = 3677.2398080981197
$ws.Range("P17").Value = 1513.1260728334094
$ws.Range("O18").Value = 1790.6868914744034
$ws.Range("P18").Value = 586.95848872727743
$ws.Range("O19").Value = 855.41816498198102
$ws.Range("P19").Value = 125.29036696573741
$ws.Range("O20").Value = 604.2894515648785
$ws.Range("P20").Value = 114.28087933685454
$ws.Range("O21").Value = 553.92428301774737
$ws.Range("P21").Value = 168.88364808995513
$ws.Range("O22").Value = 691.43918742017877
$ws.Range("P22").Value = 183.26835110673994
$ws.Range("O23").Value = 538.38801644316732
$ws.Range("P23").Value = 122.7575456246515
$ws.Range("O24").Value = 135.8047085719584
$ws.Range("P24").Value = 47.440191316611013
$ws.Range("O25").Value = 84.740911693379076
$ws.Range("P25").Value = 31.606962889645008
$ws.Range("O26").Value = 389.01613188944566
$ws.Range("P26").Value = 183.35353000496002
$ws.Range("O27").Value = 705.63269952560256
$ws.Range("P27").Value = 283.36501249227683
$ws.Range("O28").Value = 877.9498416654792
$ws.Range("P28").Value = 285.69945366645271
$ws.Range("O29").Value = 1288.3494900344381
$ws.Range("P29").Value = 353.26382814377189
$ws.Range("O30").Value = 1261.6898732778029
$ws.Range("P30").Value = 374.08152310602077
$ws.Range("O31").Value = 1320.6882719499338
$ws.Range("P31").Value = 357.81913928136248
$ws.Range("O32").Value = 1619.2040325628702
$ws.Range("P32").Value = 439.71763724873756
$ws.Range("O33").Value = 1792.1935255685628
$ws.Range("P33").Value = 382.97866840220581
$ws.Range("O34").Value = 1906.024291026303
$ws.Range("P34").Value = 379.00519413871899
$ws.Range("O35").Value = 1572.2573155346765
$ws.Range("P35").Value = 424.86006599213249
$ws.Range("O36").Value = 1528.4838728682173
$ws.Range("P36").Value = 309.94354488991985
$ws.Range("O37").Value = 1890.5007396328685
$ws.Range("P37").Value = 493.27112163403336
$ws.Range("O38").Value = 1989.2625158770416
$ws.Range("P38").Value = 701.9557933432659
$ws.Range("O39").Value = 1855.3261441497395
$ws.Range("P39").Value = 467.0670608117893
$ws.Range("O40").Value = 3940.6454600786842
$ws.Range("P40").Value = 1143.1339665169069
$ws.Range("O41").Value = 2036.2425266360976
$ws.Range("P41").Value = 354.77592678409894
$ws.Range("O42").Value = 1727.8118176156586
$ws.Range("P42").Value = 291.8278835381451
$ws.Range("O43").Value = 1474.5796351149672
$ws.Range("P43").Value = 225.57646625049392
$ws.Range("O44").Value = 1250.9863792015792
$ws.Range("P44").Value = 234.1696268731057
$ws.Range("O45").Value = 1661.3108664386739
$ws.Range("P45").Value = 349.24687261918882
$ws.Range("O46").Value = 1239.3874065421342
$ws.Range("P46").Value = 377.93930868525126
$ws.Range("O47").Value = 553.6875854764379
$ws.Range("P47").Value = 167.15378413723346
$ws.Range("O48").Value = 381.1568442259404
$ws.Range("P48").Value = 104.64998514589992
$ws.Range("O49").Value = 345.04235611416243
$ws.Range("P49").Value = 96.40136943660292
$ws.Range("O50").Value = 405.18633241906849
$ws.Range("P50").Value = 100.57567248092614
$ws.Range("O51").Value = 502.86135057239034
$ws.Range("P51").Value = 143.23164367131866
$ws.Range("O52").Value = 609.86578383430708
$ws.Range("P52").Value = 140.91601621588853
$ws.Range("O53").Value = 628.74599378661958
$ws.Range("P53").Value = 178.18833917324125
$ws.Range("O54").Value = 777.90386346391608
$ws.Range("P54").Value = 254.56241851651271
$ws.Range("O55").Value = 819.87427860569107
$ws.Range("P55").Value = 227.33965868704396
$ws.Range("O56").Value = 800.55066289636613
$ws.Range("P56").Value = 249.9419719018741
$ws.Range("O57").Value = 692.46714717198768
$ws.Range("P57").Value = 161.54193206805655
$ws.Range("O58").Value = 751.01425229442748
$ws.Range("P58").Value = 169.98534173436673
$ws.Range("O59").Value = 768.38054864868434
$ws.Range("P59").Value = 156.08475178294529
$ws.Range("O60").Value = 716.54902569990747
$ws.Range("P60").Value = 160.87595102179287
$ws.Range("O61").Value = 648.17679295783716
$ws.Range("P61").Value = 159.02501919219446
$ws.Range("O62").Value = 1759.5391597467874
$ws.Range("P62").Value = 427.54812947467633
$ws.Range("O63").Value = 964.25151149946817
$ws.Range("P63").Value = 200.6676595603054
$ws.Range("O64").Value = 675.02612769402242
$ws.Range("P64").Value = 168.70970064698326
$ws.Range("O65").Value = 521.81100256789966
$ws.Range("P65").Value = 150.53361893188935
$ws.Range("O66").Value = 435.84169482412852
$ws.Range("P66").Value = 157.17834068467118
$ws.Range("O67").Value = 689.819406004916
$ws.Range("P67").Value = 255.33802191077586
$ws.Range("O68").Value = 635.98771358848114
$ws.Range("P68").Value = 222.63633082904892
$ws.Range("O69").Value = 101.24748633857224
$ws.Range("P69").Value = 57.870559127304546
$ws.Range("O70").Value = 43.446153732682433
$ws.Range("P70").Value = 19.481130681009965
$ws.Range("O71").Value = 233.09003714234134
$ws.Range("P71").Value = 73.835653805418957
$ws.Range("O72").Value = 525.08735140180067
$ws.Range("P72").Value = 147.59747939020144
$ws.Range("O73").Value = 729.78353756724982
$ws.Range("P73").Value = 188.88548914135833
$ws.Range("O74").Value = 1269.5068328253928
$ws.Range("P74").Value = 250.69678692666332
$ws.Range("O75").Value = 1442.867264171098
$ws.Range("P75").Value = 248.35870839625946
$ws.Range("O76").Value = 1396.8787713297468
$ws.Range("P76").Value = 197.78485031247564
$ws.Range("O77").Value = 1518.7456781473618
$ws.Range("P77").Value = 173.44180797259369
$ws.Range("O78").Value = 1493.6812088852755
$ws.Range("P78").Value = 191.56076941151991
$ws.Range("O79").Value = 1416.2329313795578
$ws.Range("P79").Value = 177.05695535930957
$ws.Range("O80").Value = 1588.540998948022
$ws.Range("P80").Value = 200.41675457017672
$ws.Range("O81").Value = 1591.2319348276831
$ws.Range("P81").Value = 127.43491260971827
$ws.Range("O82").Value = 1602.2802056376374
$ws.Range("P82").Value = 142.86324795189168
$ws.Range("O83").Value = 1789.6235796883129
$ws.Range("P83").Value = 169.27996207719642
$ws.Range("O84").Value = 1728.0761897073371
$ws.Range("P84").Value = 161.71274297630765
$ws.Range("O85").Value = 2803.1554334332782
$ws.Range("P85").Value = 463.03615246140919
$ws.Range("O86").Value = 2163.5786978892343
$ws.Range("P86").Value = 213.76418640364801
$ws.Range("O87").Value = 1642.3069691780795
$ws.Range("P87").Value = 157.68446954082407
$ws.Range("O88").Value = 1270.8576160651719
$ws.Range("P88").Value = 101.73214796763017
$ws.Range("O89").Value = 1150.4877790263338
$ws.Range("P89").Value = 141.84579524582941
$ws.Range("O90").Value = 1911.6565596553658
$ws.Range("P90").Value = 265.98246709956334
$ws.Range("O91").Value = 1308.7104571620723
$ws.Range("P91").Value = 186.40924439033535

# Narrow the sem_Intake column slightly (engine quantizes ColumnWidth to
# whole internal width units, so 10.85 is the closest settable value that
# lands on the target stored width of ~11.71)
$ws.Columns.Item(16).ColumnWidth = 10.85